$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date (Excel serial number)
# for every data row. Bump it from 45171 (2023-09-02) to 45172 (2023-09-03)
# for all data rows (2 through 367).
for ($r = 2; $r -le 367; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
